$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Fix the "CDD" row (row 7): correct the misspelled mapped
#    component "Workflow Configuator" -> "Workflow Configurator", and
#    give the row the thick-bottom "separator" styling that matches
#    the last row's look (row height 15, thick bottom border).
# ------------------------------------------------------------------
$ws.Range("B7").Value2 = "Workflow Configurator"
$ws.Rows.Item(7).RowHeight = 15

$ws.Range("B7").Borders.Item(7).LineStyle = 1
$ws.Range("B7").Borders.Item(7).Weight = 4
$ws.Range("B7").Borders.Item(10).LineStyle = 1
$ws.Range("B7").Borders.Item(10).Weight = 2
$ws.Range("B7").Borders.Item(8).LineStyle = 1
$ws.Range("B7").Borders.Item(8).Weight = 2
$ws.Range("B7").Borders.Item(9).LineStyle = 1
$ws.Range("B7").Borders.Item(9).Weight = 4
$ws.Range("B7").Font.Color = $ws.Range("B4").Font.Color

# ------------------------------------------------------------------
# 2) Insert a new table row above row 13 ("Other" / "Scenario
#    Manager" boundary) for the new "P&C" component, keeping the
#    table sorted alphabetically by Component.
# ------------------------------------------------------------------
$ws.Range("A13:C13").Insert(-4121)

$ws.Range("A5:C5").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)
$ws.Range("A13:B13").Font.Name = "Calibri"
$ws.Range("C13").Font.Name = "Calibri"

$ws.Range("A13").Value2 = "P&C "
$ws.Range("B13").Value2 = "Base Platform"
$ws.Range("C13").Value2 = "England"

# Keep the table definition in sync with the newly inserted row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C17"))

# ------------------------------------------------------------------
# 3) Leave the selection where the edit finished, matching the
#    post-edit cursor position recorded in the workbook.
# ------------------------------------------------------------------
$ws.Range("E4").Select()
